$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Refresh the crypto price ("Price", column D) and volume-change ("Volume(1h)", column E) data ---
# Column D mixes locale-formatted numbers (dot-separated thousands, e.g. "30.321.93") with plain
# decimals (e.g. "345.39"). Force the whole price column to Text format first so Excel keeps every
# new value exactly as scraped instead of auto-parsing it into a numeric/date value; then restore
# the original (default) formatting once all values are written.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '30.321.93'
$ws.Range("E2").Value = '  +2.27%  '
$ws.Range("D3").Value = '2.106.49'
$ws.Range("E3").Value = '  +0.81%  '
$ws.Range("D4").Value = '1.006'
$ws.Range("E4").Value = '  -0.30%  '
$ws.Range("D5").Value = '345.39'
$ws.Range("E5").Value = '  +0.68%  '
$ws.Range("D6").Value = '1.006'
$ws.Range("E6").Value = '  -0.19%  '
$ws.Range("D7").Value = '0.5225'
$ws.Range("E7").Value = '  +1.35%  '
$ws.Range("D8").Value = '0.4440'
$ws.Range("E8").Value = '  +1.11%  '
$ws.Range("D9").Value = '54.48'
$ws.Range("E9").Value = '  +4.49%  '
$ws.Range("D10").Value = '0.09447'
$ws.Range("E10").Value = '  +2.17%  '
$ws.Range("E11").Value = '  -0.25%  '
$ws.Range("D12").Value = '25.12'
$ws.Range("E12").Value = '  +0.23%  '
$ws.Range("D13").Value = '8.724'
$ws.Range("E13").Value = '  +6.88%  '
$ws.Range("D14").Value = '2.141.61'
$ws.Range("E14").Value = '  +2.43%  '
$ws.Range("D15").Value = '6.920'
$ws.Range("E15").Value = '  +2.41%  '
$ws.Range("D16").Value = '101.79'
$ws.Range("E16").Value = '  +1.67%  '
$ws.Range("D17").Value = '0.00001163'
$ws.Range("E17").Value = '  +0.67%  '
$ws.Range("D18").Value = '1.007'
$ws.Range("E18").Value = '  -0.19%  '
$ws.Range("D19").Value = '21.33'
$ws.Range("E19").Value = '  +1.22%  '
$ws.Range("D20").Value = '0.06730'
$ws.Range("E20").Value = '  +1.25%  '
$ws.Range("D21").Value = '6.317'
$ws.Range("E21").Value = '  +2.21%  '
$ws.Range("D22").Value = '1.006'
$ws.Range("E22").Value = '  -0.21%  '
$ws.Range("D23").Value = '30.359.84'
$ws.Range("E23").Value = '  +2.25%  '
$ws.Range("D24").Value = '12.64'
$ws.Range("E24").Value = '  +0.02%  '
$ws.Range("D25").Value = '2.317'
$ws.Range("E25").Value = '  +0.06%  '
$ws.Range("D26").Value = '2.370.84'
$ws.Range("E26").Value = '  +1.43%  '
$ws.Range("D27").Value = '22.03'
$ws.Range("E27").Value = '  +0.75%  '
$ws.Range("D28").Value = '2.545'
$ws.Range("E28").Value = '  +0.65%  '
$ws.Range("D29").Value = '163.06'
$ws.Range("E29").Value = '  +0.31%  '
$ws.Range("E30").Value = '  +0.66%  '
$ws.Range("D31").Value = '1.147'
$ws.Range("E31").Value = '  +0.97%  '
$ws.Range("D32").Value = '1.766'
$ws.Range("E32").Value = '  +8.58%  '
$ws.Range("D33").Value = '0.1056'
$ws.Range("E33").Value = '  +0.77%  '
$ws.Range("D34").Value = '6.856'
$ws.Range("E34").Value = '  +12.92%  '
$ws.Range("D35").Value = '6.260'
$ws.Range("E35").Value = '  +1.27%  '
$ws.Range("D36").Value = '3.921'
$ws.Range("E36").Value = '  -1.02%  '
$ws.Range("D37").Value = '10.55'
$ws.Range("E37").Value = '  +2.80%  '
$ws.Range("D38").Value = '0.02632'
$ws.Range("E38").Value = '  +2.52%  '
$ws.Range("D39").Value = '0.06800'
$ws.Range("E39").Value = '  +1.18%  '
$ws.Range("D40").Value = '0.7045'
$ws.Range("E40").Value = '  +2.05%  '
$ws.Range("D41").Value = '12.56'
$ws.Range("E41").Value = '  +1.30%  '
$ws.Range("D42").Value = '1.344'
$ws.Range("E42").Value = '  +4.28%  '
$ws.Range("E43").Value = '  -0.37%  '
$ws.Range("E44").Value = '  +1.29%  '
$ws.Range("D45").Value = '14.50'
$ws.Range("E45").Value = '  +3.00%  '
$ws.Range("D46").Value = '2.363'
$ws.Range("E46").Value = '  +2.43%  '
$ws.Range("D47").Value = '1.005'
$ws.Range("E47").Value = '  -0.18%  '
$ws.Range("D48").Value = '1.363'
$ws.Range("E48").Value = '  +16.94%  '
$ws.Range("D49").Value = '3.652'
$ws.Range("E49").Value = '  +1.11%  '
$ws.Range("D50").Value = '0.00000000344'
$ws.Range("E50").Value = '  +0.96%  '
$ws.Range("D51").Value = '1.222'
$ws.Range("E51").Value = '  +0.41%  '

# Restore the default (unformatted) style on column D now that the text values are in place.
$priceRange.ClearFormats()
